$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (autogluon) - previously empty inline strings, now populated
$ws.Range("B3").Value = "0.399 (0.350 ± 0.021)"
$ws.Range("C3").Value = "00:02:39 (00:02:48 ± 00:00:07)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("E3").Value = "[]"

# F3 looks like a plain integer, so force text storage (matches sibling
# cells such as F4/F6/F8 which are stored as text, not numbers), then
# strip the number-format styling that forcing text applied so the cell
# keeps the sheet's default (unstyled) appearance.
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "19"
$ws.Range("F3").ClearFormats()

# Row 4 (autokeras) - fix mojibake "Â±" -> "±"
$ws.Range("B4").Value = "0.711 (0.677 ± 0.016)"
$ws.Range("C4").Value = "00:03:15 (00:03:49 ± 00:00:36)"
$ws.Range("D4").Value = "00:00:10 (00:00:10 ± 00:00:00)"

# Row 6 (autosklearn) - fix mojibake "Â±" -> "±"
$ws.Range("B6").Value = "0.807 (0.775 ± 0.015)"
$ws.Range("C6").Value = "00:04:56 (00:05:01 ± 00:00:02)"
$ws.Range("D6").Value = "00:00:00 (00:00:02 ± 00:00:01)"

# Row 8 (fedot) - fix mojibake "Â±" -> "±"
$ws.Range("B8").Value = "0.744 (0.689 ± 0.030)"
$ws.Range("C8").Value = "00:05:06 (00:09:55 ± 00:04:20)"
$ws.Range("D8").Value = "00:00:00 (00:00:00 ± 00:00:00)"
